$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 357
$ws.Range("F5").Value = 40
$ws.Range("F6").Value = 6131
$ws.Range("F7").Value = 677
$ws.Range("F8").Value = 1071
$ws.Range("F9").Value = 33
$ws.Range("F10").Value = 280
$ws.Range("F11").Value = 184
$ws.Range("F13").Value = 616
$ws.Range("F14").Value = 1057
$ws.Range("F15").Value = 62
$ws.Range("F17").Value = 296
$ws.Range("F18").Value = 1380
$ws.Range("F20").Value = 1035
$ws.Range("F21").Value = 85
$ws.Range("F22").Value = 2102
$ws.Range("F23").Value = 210
$ws.Range("F24").Value = 53
$ws.Range("F27").Value = 3405

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 26
$ws.Range("F5").Value = 110
$ws.Range("F8").Value = 22
$ws.Range("F9").Value = 670
$ws.Range("F14").Value = 87
$ws.Range("F20").Value = 306
$ws.Range("F21").Value = 4066
$ws.Range("F22").Value = 10
$ws.Range("F24").Value = 29
$ws.Range("F25").Value = 165
$ws.Range("F27").Value = 79

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1808
$ws.Range("F6").Value = 1163
$ws.Range("F8").Value = 1527
$ws.Range("F12").Value = 697

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1808
$ws.Range("F5").Value = 1163
$ws.Range("F6").Value = 1527
$ws.Range("F9").Value = 697
$ws.Range("F12").Value = 357
$ws.Range("F13").Value = 40
$ws.Range("F14").Value = 6131
$ws.Range("F15").Value = 22
$ws.Range("F16").Value = 677
$ws.Range("F17").Value = 1071
$ws.Range("F18").Value = 33
$ws.Range("F19").Value = 280
$ws.Range("F20").Value = 184
$ws.Range("F22").Value = 616
$ws.Range("F25").Value = 87
$ws.Range("F28").Value = 1057
$ws.Range("F29").Value = 62
$ws.Range("F31").Value = 296
$ws.Range("F34").Value = 1380
$ws.Range("F36").Value = 29
$ws.Range("F37").Value = 165
$ws.Range("F39").Value = 1035
$ws.Range("F40").Value = 85
$ws.Range("F42").Value = 2102
$ws.Range("F44").Value = 210
$ws.Range("F45").Value = 53
$ws.Range("F48").Value = 3405
